# Refresh the cryptocurrency price/volume table with the latest scraped values
# (mirrors the automated "Updated cryptos list ... with GitHub Actions" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # These Price/Volume columns store plain text (e.g. "26.888.85", "0.9978",
    # "  +1.81%  ") rather than real numbers. Force the cell to Text format so
    # Excel does not auto-coerce numeric-looking input into a number, then restore
    # the cell to its original (unstyled) state so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '26.888.85'
Set-TextValue $ws.Range("E2") '  +1.81%  '
Set-TextValue $ws.Range("D3") '1.726.46'
Set-TextValue $ws.Range("E3") '  +0.22%  '
Set-TextValue $ws.Range("E4") '  -0.24%  '
Set-TextValue $ws.Range("D5") '241.94'
Set-TextValue $ws.Range("E5") '  -0.34%  '
Set-TextValue $ws.Range("D6") '0.9978'
Set-TextValue $ws.Range("E6") '  -0.21%  '
Set-TextValue $ws.Range("D7") '0.4889'
Set-TextValue $ws.Range("E7") '  -0.70%  '
Set-TextValue $ws.Range("D8") '0.2593'
Set-TextValue $ws.Range("E8") '  -0.86%  '
Set-TextValue $ws.Range("D9") '0.06211'
Set-TextValue $ws.Range("E9") '  +0.21%  '
Set-TextValue $ws.Range("D10") '1.731.51'
Set-TextValue $ws.Range("E10") '  +0.53%  '
Set-TextValue $ws.Range("D11") '15.99'
Set-TextValue $ws.Range("E11") '  +3.21%  '
Set-TextValue $ws.Range("D12") '0.06901'
Set-TextValue $ws.Range("E12") '  -1.67%  '
Set-TextValue $ws.Range("D13") '0.6080'
Set-TextValue $ws.Range("E13") '  +1.42%  '
Set-TextValue $ws.Range("D14") '4.484'
Set-TextValue $ws.Range("E14") '  -1.93%  '
Set-TextValue $ws.Range("D15") '77.21'
Set-TextValue $ws.Range("E15") '  -0.02%  '
Set-TextValue $ws.Range("E16") '  -0.18%  '
Set-TextValue $ws.Range("D17") '26.645.75'
Set-TextValue $ws.Range("E17") '  +0.91%  '
Set-TextValue $ws.Range("D18") '0.9973'
Set-TextValue $ws.Range("E18") '  -0.25%  '
Set-TextValue $ws.Range("D19") '0.000007177'
Set-TextValue $ws.Range("E19") '  +0.12%  '
Set-TextValue $ws.Range("D20") '11.44'
Set-TextValue $ws.Range("E20") '  +0.79%  '
Set-TextValue $ws.Range("D21") '1.953.71'
Set-TextValue $ws.Range("E21") '  +0.42%  '
Set-TextValue $ws.Range("D22") '4.424'
Set-TextValue $ws.Range("E22") '  -1.31%  '
Set-TextValue $ws.Range("D23") '8.562'
Set-TextValue $ws.Range("E23") '  -0.25%  '
Set-TextValue $ws.Range("D24") '5.097'
Set-TextValue $ws.Range("E24") '  -1.22%  '
Set-TextValue $ws.Range("D25") '138.40'
Set-TextValue $ws.Range("E25") '  +0.76%  '
Set-TextValue $ws.Range("D26") '15.31'
Set-TextValue $ws.Range("E26") '  +0.53%  '
Set-TextValue $ws.Range("D27") '1.776'
Set-TextValue $ws.Range("E27") '  +4.10%  '
Set-TextValue $ws.Range("D28") '106.37'
Set-TextValue $ws.Range("E28") '  -0.66%  '
Set-TextValue $ws.Range("D29") '1.381'
Set-TextValue $ws.Range("E29") '  -1.16%  '
Set-TextValue $ws.Range("D30") '3.942'
Set-TextValue $ws.Range("E30") '  +0.12%  '
Set-TextValue $ws.Range("D31") '0.08002'
Set-TextValue $ws.Range("E31") '  +0.46%  '
Set-TextValue $ws.Range("E33") '  -0.40%  '
Set-TextValue $ws.Range("E34") '  -0.29%  '
Set-TextValue $ws.Range("E35") '  +1.40%  '
Set-TextValue $ws.Range("D36") '0.6253'
Set-TextValue $ws.Range("E36") '  +0.01%  '
Set-TextValue $ws.Range("D37") '0.9369'
Set-TextValue $ws.Range("E37") '  +1.43%  '
Set-TextValue $ws.Range("D38") '2.051'
Set-TextValue $ws.Range("E38") '  +5.31%  '
Set-TextValue $ws.Range("D39") '2.452'
Set-TextValue $ws.Range("E39") '  +2.52%  '
Set-TextValue $ws.Range("D40") '0.9974'
Set-TextValue $ws.Range("E40") '  -0.21%  '
Set-TextValue $ws.Range("B41") 'FraxShare'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D41") '5.701'
Set-TextValue $ws.Range("E41") '  +7.02%  '
Set-TextValue $ws.Range("B42") 'VeChain'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D42") '0.01501'
Set-TextValue $ws.Range("E42") '  +1.21%  '
Set-TextValue $ws.Range("D43") '99.61'
Set-TextValue $ws.Range("E43") '  -0.29%  '
Set-TextValue $ws.Range("E44") '  +0.35%  '
Set-TextValue $ws.Range("D45") '6.875'
Set-TextValue $ws.Range("E45") '  +2.19%  '
Set-TextValue $ws.Range("E46") '  -0.23%  '
Set-TextValue $ws.Range("D47") '0.05401'
Set-TextValue $ws.Range("E47") '  +0.68%  '
Set-TextValue $ws.Range("D48") '7.942'
Set-TextValue $ws.Range("E48") '  +3.52%  '
Set-TextValue $ws.Range("D49") '30.18'
Set-TextValue $ws.Range("E49") '  +0.33%  '
Set-TextValue $ws.Range("E50") '  +1.60%  '
Set-TextValue $ws.Range("D51") '1.233'
Set-TextValue $ws.Range("E51") '  -0.05%  '
